$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# low_var_feature_removal: turn on (1) and request top 15 features -> "1, 15"
$ws.Range("B17").Value = "1, 15"

# Update the active selection on the sheet to A18:C18 (the
# "database injection settings" section header row)
$ws.Range("A18:C18").Select()
